$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowValues = @(2026, 2, 11, 5201, 4959, 4886, 4856, 4875, 5049, 5476, 5985, 6402, 6570, 6465, 6264, 6069, 5909, 5819, 5879, 5899, 6167, 6545, 6535, 6318, 6024, 5752, 5395)

$row = 8
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item($row, $col).Value = $rowValues[$i]
}
